$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D3").Value = -7.977000000000001
$ws.Range("C7").Value = -13.682
$ws.Range("B10").Value = 6.415000000000001
$ws.Range("E10").Value = 16.443
$ws.Range("B12").Value = 5.86
$ws.Range("E14").Value = 17.027
$ws.Range("C15").Value = -13.741
$ws.Range("B18").Value = 5.249
$ws.Range("D18").Value = -8.784000000000001
$ws.Range("D19").Value = -8.113000000000001
$ws.Range("C20").Value = -12.673
$ws.Range("D27").Value = -8.624000000000001
$ws.Range("C29").Value = -12.24
$ws.Range("C30").Value = -12.981
$ws.Range("C31").Value = -13.105
$ws.Range("E32").Value = 17.041
$ws.Range("E35").Value = 16.358
$ws.Range("B37").Value = 8.494
$ws.Range("C40").Value = -12.782
$ws.Range("D42").Value = -8.468
$ws.Range("E43").Value = 17.212
$ws.Range("D44").Value = -7.836
$ws.Range("D47").Value = -7.605999999999999
$ws.Range("E49").Value = 16.349
$ws.Range("B55").Value = 5.194
$ws.Range("E56").Value = 16.166
$ws.Range("D58").Value = -8.403
$ws.Range("B68").Value = 5.403
$ws.Range("C68").Value = -11.083
$ws.Range("E69").Value = 17.438
$ws.Range("D73").Value = -7.841999999999999
$ws.Range("C76").Value = -12.977
$ws.Range("B77").Value = 6.353
$ws.Range("B78").Value = 7.631
$ws.Range("E81").Value = 16.583
$ws.Range("C87").Value = -13.393
$ws.Range("C88").Value = -13.476
$ws.Range("E92").Value = 17.853
$ws.Range("D95").Value = -7.564000000000002
$ws.Range("C96").Value = -12.628
$ws.Range("C98").Value = -13.294
$ws.Range("C101").Value = -12.612
$ws.Range("D101").Value = -7.768000000000001
$ws.Range("C102").Value = -13.086
